# Monte Carlo simulation setup: populate the "stdev" (column C) values for
# the environmental config parameters, and widen column C to fit (matching
# column B) since it now holds numbers of similar width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Existing stdev values that were tightened/re-estimated.
$ws.Range("C31").Value = 0.01
$ws.Range("C32").Value = 0.01
$ws.Range("C33").Value = 0.01
$ws.Range("C34").Value = 0.01

# Previously-blank stdev cells now filled in with Monte Carlo inputs.
$ws.Range("C35").Value = 0.01
$ws.Range("C36").Value = 0.01
$ws.Range("C37").Value = 0.01
$ws.Range("C38").Value = 0.01
$ws.Range("C39").Value = 0.02
$ws.Range("C40").Value = 0.01
$ws.Range("C41").Value = 0.04
$ws.Range("C42").Value = 0.002
$ws.Range("C43").Value = 0.01
$ws.Range("C44").Value = 0.02
$ws.Range("C45").Value = 0.05
$ws.Range("C46").Value = 150
$ws.Range("C52").Value = 0.000005
$ws.Range("C53").Value = 0.0001

$ws.Range("C66").Value = 0.01
$ws.Range("C67").Value = 0.01
$ws.Range("C68").Value = 0.01
$ws.Range("C69").Value = 0.04

# Column C now carries values comparable in width to column B, so fit it.
$ws.Columns("C").AutoFit() | Out-Null
